# "aggiornamento fino a 28 luglio" -- append the next batch of daily rows
# (new-positive / 7-day rolling sum / 7-day rolling sum per 100k) to the
# bottom of the single data sheet, continuing the existing A:D table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 301
$firstNewRow = $lastRow + 1

# date serial (col A), nuovi pos. (col B), somma mobile 7gg. (col C),
# somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
    @(44376, 0, 0, 0),
    @(44377, 0, 0, 0),
    @(44378, 0, 0, 0),
    @(44379, 0, 0, 0),
    @(44380, 0, 0, 0),
    @(44381, 0, 0, 0),
    @(44382, 0, 0, 0),
    @(44383, 0, 0, 0),
    @(44384, 0, 0, 0),
    @(44385, 0, 0, 0),
    @(44386, 0, 0, 0),
    @(44387, 2, 2, 17.71636105943839),
    @(44388, 0, 2, 17.71636105943839),
    @(44389, 0, 2, 17.71636105943839),
    @(44390, 0, 2, 17.71636105943839),
    @(44391, 0, 2, 17.71636105943839),
    @(44392, 0, 2, 17.71636105943839),
    @(44393, 0, 2, 17.71636105943839),
    @(44394, 0, 0, 0),
    @(44395, 0, 0, 0),
    @(44396, 0, 0, 0),
    @(44397, 0, 0, 0),
    @(44398, 0, 0, 0),
    @(44399, 0, 0, 0),
    @(44400, 2, 2, 17.71636105943839),
    @(44401, 3, 5, 44.29090264859597),
    @(44402, 0, 5, 44.29090264859597)
)

$lastNewRow = $firstNewRow + $data.Length - 1

# Extend column A's date formatting/style down the new rows by copying the
# format from the current last row (keeps the same cell style as the rest
# of the date column) before writing values.
$ws.Range("A$lastRow").Copy() | Out-Null
$ws.Range("A$firstNewRow`:A$lastNewRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $firstNewRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
}
